# Auto update stock data
# Updates the "Date_1" (A2) and "EBITDA" (B2) figures on each company sheet to
# the latest daily refresh, plus the extra Ryerson Holding (RYI) adjustments
# (EBITDA, Inventory Turnover, Current Ratio updates and clearing the stale
# Altman Z-Score column for RYI rows 2-8).

$wb = $excel.ActiveWorkbook

# Helper: write a value to a cell while forcing Text storage so that
# numeric/date-looking strings ("2025/10/29", "4.78", ...) are kept as text
# instead of being auto-converted by Excel into a date serial / number.
function Set-TextValue($ws, $cellRef, $text) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $text
}

# --- Sheet: Alcoa (AA) ---
$ws = $wb.Worksheets.Item("Alcoa")
Set-TextValue $ws "A2" "2025/10/29"
Set-TextValue $ws "B2" "4.78"

# --- Sheet: Rio Tinto (RIO) ---
$ws = $wb.Worksheets.Item("Rio Tinto")
Set-TextValue $ws "A2" "2025/10/29"
Set-TextValue $ws "B2" "7.69"

# --- Sheet: Norsk Hydro (NHY) ---
$ws = $wb.Worksheets.Item("Norsk Hydro")
Set-TextValue $ws "A2" "2025/10/29"
Set-TextValue $ws "B2" "2.70"

# --- Sheet: Reliance Steel & Aluminum (RS) ---
$ws = $wb.Worksheets.Item("Reliance Steel & Aluminum")
Set-TextValue $ws "A2" "2025/10/29"
Set-TextValue $ws "B2" "12.39"

# --- Sheet: Kaiser Aluminum (KALU) ---
$ws = $wb.Worksheets.Item("Kaiser Aluminum")
Set-TextValue $ws "A2" "2025/10/29"
Set-TextValue $ws "B2" "9.82"

# --- Sheet: Ryerson Holding (RYI) ---
$ws = $wb.Worksheets.Item("Ryerson Holding")
Set-TextValue $ws "A2" "2025/10/29"
Set-TextValue $ws "B2" "26.43"
Set-TextValue $ws "D2" "5.40"
Set-TextValue $ws "E2" "1.98"

# The Altman Z-Score column (G) for Ryerson Holding is no longer populated
# for any of the historical rows - clear it out.
$ws.Range("G2:G8").ClearContents()
